$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.644.68"
$ws.Range("E2").Value = "'  +2.38%  "
$ws.Range("D3").Value = "'3.813.73"
$ws.Range("E3").Value = "'  +1.28%  "
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D5").Value = "'664.55"
$ws.Range("E5").Value = "'  +7.08%  "
$ws.Range("D6").Value = "'168.85"
$ws.Range("E6").Value = "'  +2.92%  "
$ws.Range("D7").Value = "'3.811.46"
$ws.Range("E7").Value = "'  +1.31%  "
$ws.Range("E8").Value = "'  -0.06%  "
$ws.Range("E9").Value = "'  +1.66%  "
$ws.Range("E10").Value = "'  +0.60%  "
$ws.Range("D12").Value = "'7.02"
$ws.Range("E12").Value = "'  +6.17%  "
$ws.Range("D13").Value = "'0.0000244"
$ws.Range("E13").Value = "'  -0.97%  "
$ws.Range("D14").Value = "'35.89"
$ws.Range("E14").Value = "'  +1.65%  "
$ws.Range("D15").Value = "'4.459.33"
$ws.Range("E15").Value = "'  +1.27%  "
$ws.Range("D16").Value = "'3.814.19"
$ws.Range("E16").Value = "'  +2.77%  "
$ws.Range("D17").Value = "'70.662.63"
$ws.Range("D18").Value = "'17.77"
$ws.Range("E18").Value = "'  +0.66%  "
$ws.Range("D19").Value = "'7.17"
$ws.Range("E19").Value = "'  +1.28%  "
$ws.Range("E20").Value = "'  +0.91%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'10.35"
$ws.Range("E21").Value = "'  +8.44%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'478.28"
$ws.Range("E22").Value = "'  +2.46%  "
$ws.Range("D23").Value = "'0.712"
$ws.Range("E23").Value = "'  +1.91%  "
$ws.Range("D24").Value = "'0.0000145"
$ws.Range("E24").Value = "'  -2.25%  "
$ws.Range("D25").Value = "'82.84"
$ws.Range("E25").Value = "'  -0.19%  "
$ws.Range("D26").Value = "'12.27"
$ws.Range("D27").Value = "'10.38"
$ws.Range("E27").Value = "'  +3.94%  "
$ws.Range("D28").Value = "'2.12"
$ws.Range("E28").Value = "'  -0.94%  "
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("D30").Value = "'3.966.47"
$ws.Range("E30").Value = "'  +1.35%  "
$ws.Range("D31").Value = "'2.85"
$ws.Range("E31").Value = "'  +7.55%  "
$ws.Range("D32").Value = "'2.31"
$ws.Range("E32").Value = "'  +3.14%  "
$ws.Range("D33").Value = "'7.40"
$ws.Range("E33").Value = "'  +1.94%  "
$ws.Range("D34").Value = "'29.54"
$ws.Range("E34").Value = "'  +2.57%  "
$ws.Range("D35").Value = "'0.180"
$ws.Range("E35").Value = "'  +14.77%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "'9.15"
$ws.Range("E36").Value = "'  +2.28%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "'  +0.04%  "
$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "'3.771.18"
$ws.Range("E38").Value = "'  +1.42%  "
$ws.Range("D40").Value = "'3.41"
$ws.Range("E40").Value = "'  +1.80%  "
$ws.Range("E41").Value = "'  +2.66%  "
$ws.Range("D42").Value = "'0.971"
$ws.Range("E42").Value = "'  +0.95%  "
$ws.Range("E43").Value = "'  +0.13%  "
$ws.Range("E44").Value = "'  -0.02%  "
$ws.Range("D45").Value = "'2.07"
$ws.Range("E45").Value = "'  +9.15%  "
$ws.Range("D46").Value = "'45.43"
$ws.Range("E46").Value = "'  +6.56%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'48.80"
$ws.Range("E47").Value = "'  +4.65%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'158.66"
$ws.Range("E48").Value = "'  +2.97%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").Value = "'1.43"
$ws.Range("E49").Value = "'  +5.12%  "
$ws.Range("D50").Value = "'0.299"
$ws.Range("E50").Value = "'  +0.42%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'8.49"
$ws.Range("E51").Value = "'  +1.28%  "
